$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2121376991271973
$ws.Range("E2").Value = 126.6717933621676
$ws.Range("F2").Value = 0.00528233382665625
$ws.Range("G2").Value = 0.00361370371671504
$ws.Range("H2").Value = 0.003382748082890474
$ws.Range("I2").Value = 0.003246370792617207
$ws.Range("J2").Value = 0.003063485040150263
$ws.Range("K2").Value = 0.003063485040150263
$ws.Range("L2").Value = 0.003061190610580106
$ws.Range("M2").Value = 0.002838050800030145
$ws.Range("N2").Value = 0.00277996133201459
$ws.Range("O2").Value = 0.002714580956402124
$ws.Range("P2").Value = 0.002700471993269063
$ws.Range("Q2").Value = 0.002700471993269063
$ws.Range("R2").Value = 0.00268194070224144
$ws.Range("S2").Value = 0.002619538620979761
$ws.Range("T2").Value = 0.002593738119387984
$ws.Range("U2").Value = 0.002579888395264869
$ws.Range("V2").Value = 0.002528398980460917
$ws.Range("W2").Value = 0.002519106432429755
$ws.Range("X2").Value = 0.002475901483295411
$ws.Range("Y2").Value = 0.002469235738053949
$ws.Range("C3").Value = 0.2116734981536865
$ws.Range("E3").Value = 131.820851839464
$ws.Range("F3").Value = 0.005332121010410525
$ws.Range("G3").Value = 0.003696987142443725
$ws.Range("H3").Value = 0.003123465910441338
$ws.Range("I3").Value = 0.003123465910441338
$ws.Range("J3").Value = 0.003123465910441338
$ws.Range("K3").Value = 0.003123465910441338
$ws.Range("L3").Value = 0.002898282955195798
$ws.Range("M3").Value = 0.002898282955195798
$ws.Range("N3").Value = 0.002898282955195798
$ws.Range("O3").Value = 0.002898282955195798
$ws.Range("P3").Value = 0.002898282955195798
$ws.Range("Q3").Value = 0.002746105463780945
$ws.Range("R3").Value = 0.002746105463780945
$ws.Range("S3").Value = 0.002699121978011494
$ws.Range("T3").Value = 0.002670520771063845
$ws.Range("U3").Value = 0.002670520771063845
$ws.Range("V3").Value = 0.002649091186681688
$ws.Range("W3").Value = 0.002640632792734297
$ws.Range("X3").Value = 0.002607742563873597
$ws.Range("Y3").Value = 0.002569607248332631
$ws.Range("C4").Value = 0.2542331218719482
$ws.Range("E4").Value = 121.2549160301805
$ws.Range("F4").Value = 0.005076595900625629
$ws.Range("G4").Value = 0.004253683968261391
$ws.Range("H4").Value = 0.00393359681751459
$ws.Range("I4").Value = 0.003564457652490048
$ws.Range("J4").Value = 0.003345836869803324
$ws.Range("K4").Value = 0.003057499272538612
$ws.Range("L4").Value = 0.00300489472005473
$ws.Range("M4").Value = 0.002926809483688932
$ws.Range("N4").Value = 0.002852985676140571
$ws.Range("O4").Value = 0.002710994630188901
$ws.Range("P4").Value = 0.002710994630188901
$ws.Range("Q4").Value = 0.002574238505891208
$ws.Range("R4").Value = 0.002520746780833594
$ws.Range("S4").Value = 0.002508689487193101
$ws.Range("T4").Value = 0.002452627223673344
$ws.Range("U4").Value = 0.002448012654612502
$ws.Range("V4").Value = 0.002424903538945128
$ws.Range("W4").Value = 0.002393421319376575
$ws.Range("X4").Value = 0.002384404872589559
$ws.Range("Y4").Value = 0.002363643587332953
$ws.Range("C5").Value = 0.2052321434020996
$ws.Range("E5").Value = 126.8747968662647
$ws.Range("F5").Value = 0.005326333466154747
$ws.Range("G5").Value = 0.004100683574910057
$ws.Range("H5").Value = 0.003692601476026889
$ws.Range("I5").Value = 0.00364151815449644
$ws.Range("J5").Value = 0.003431845201106533
$ws.Range("K5").Value = 0.003292767639728839
$ws.Range("L5").Value = 0.003292767639728839
$ws.Range("M5").Value = 0.003145769079246658
$ws.Range("N5").Value = 0.00305566689409799
$ws.Range("O5").Value = 0.002911001091240046
$ws.Range("P5").Value = 0.0028794657918634
$ws.Range("Q5").Value = 0.002860307584064622
$ws.Range("R5").Value = 0.002814242707251737
$ws.Range("S5").Value = 0.002700703172838516
$ws.Range("T5").Value = 0.002581676515431161
$ws.Range("U5").Value = 0.002581676515431161
$ws.Range("V5").Value = 0.002533499274059582
$ws.Range("W5").Value = 0.002508942230674646
$ws.Range("X5").Value = 0.002491968678754928
$ws.Range("Y5").Value = 0.002473192921369681
$ws.Range("C6").Value = 0.2187762260437012
$ws.Range("E6").Value = 129.0558737381361
$ws.Range("F6").Value = 0.004973533463168622
$ws.Range("G6").Value = 0.004007834176334806
$ws.Range("H6").Value = 0.003765545263563394
$ws.Range("I6").Value = 0.00329902033455459
$ws.Range("J6").Value = 0.002996318119233419
$ws.Range("K6").Value = 0.002837164056632992
$ws.Range("L6").Value = 0.002837164056632992
$ws.Range("M6").Value = 0.002754386194236799
$ws.Range("N6").Value = 0.002617625099652148
$ws.Range("O6").Value = 0.002617625099652148
$ws.Range("P6").Value = 0.002617625099652148
$ws.Range("Q6").Value = 0.002617625099652148
$ws.Range("R6").Value = 0.002617625099652148
$ws.Range("S6").Value = 0.002595882651043438
$ws.Range("T6").Value = 0.00259323266272408
$ws.Range("U6").Value = 0.002579049413439581
$ws.Range("V6").Value = 0.002551534453327475
$ws.Range("W6").Value = 0.002531023254919915
$ws.Range("X6").Value = 0.002526609723334271
$ws.Range("Y6").Value = 0.002515709039729749
$ws.Range("C7").Value = 0.2400655746459961
$ws.Range("E7").Value = 128.0566292336225
$ws.Range("F7").Value = 0.005219006066630359
$ws.Range("G7").Value = 0.004051474292314311
$ws.Range("H7").Value = 0.003906398998612191
$ws.Range("I7").Value = 0.003188090128879357
$ws.Range("J7").Value = 0.003128614044360989
$ws.Range("K7").Value = 0.003119278016120661
$ws.Range("L7").Value = 0.002904662338371752
$ws.Range("M7").Value = 0.002904662338371752
$ws.Range("N7").Value = 0.002904662338371752
$ws.Range("O7").Value = 0.002904662338371752
$ws.Range("P7").Value = 0.002732313833185527
$ws.Range("Q7").Value = 0.002732313833185527
$ws.Range("R7").Value = 0.002665358665284617
$ws.Range("S7").Value = 0.002658306804732575
$ws.Range("T7").Value = 0.002634189869973749
$ws.Range("U7").Value = 0.00256335318610163
$ws.Range("V7").Value = 0.00256335318610163
$ws.Range("W7").Value = 0.002529427871117678
$ws.Range("X7").Value = 0.002510368390949313
$ws.Range("Y7").Value = 0.002496230589349366
$ws.Range("C8").Value = 0.1998984813690186
$ws.Range("E8").Value = 123.7813224347901
$ws.Range("F8").Value = 0.004942089665453015
$ws.Range("G8").Value = 0.004137135732418034
$ws.Range("H8").Value = 0.003509129963556127
$ws.Range("I8").Value = 0.003398433446092955
$ws.Range("J8").Value = 0.003268118359936312
$ws.Range("K8").Value = 0.003139955469954653
$ws.Range("L8").Value = 0.002907431792588261
$ws.Range("M8").Value = 0.002907431792588261
$ws.Range("N8").Value = 0.002809764747183962
$ws.Range("O8").Value = 0.002759037451199575
$ws.Range("P8").Value = 0.002642192144686069
$ws.Range("Q8").Value = 0.002585381755691459
$ws.Range("R8").Value = 0.002585381755691459
$ws.Range("S8").Value = 0.002564262630548981
$ws.Range("T8").Value = 0.002517860241060849
$ws.Range("U8").Value = 0.002508724663945088
$ws.Range("V8").Value = 0.002495515038345404
$ws.Range("W8").Value = 0.002469520495740167
$ws.Range("X8").Value = 0.002424979380251238
$ws.Range("Y8").Value = 0.002412891275531971
$ws.Range("C9").Value = 0.2187156677246094
$ws.Range("E9").Value = 123.1594028764302
$ws.Range("F9").Value = 0.005051936418252378
$ws.Range("G9").Value = 0.003993786927767359
$ws.Range("H9").Value = 0.003485043839811612
$ws.Range("I9").Value = 0.003386969770129717
$ws.Range("J9").Value = 0.003262911120083147
$ws.Range("K9").Value = 0.002897553220910092
$ws.Range("L9").Value = 0.002823757717874526
$ws.Range("M9").Value = 0.002808924701578921
$ws.Range("N9").Value = 0.002808924701578921
$ws.Range("O9").Value = 0.002710307137762388
$ws.Range("P9").Value = 0.00265108004136262
$ws.Range("Q9").Value = 0.002614194224425041
$ws.Range("R9").Value = 0.002614194224425041
$ws.Range("S9").Value = 0.00254398382717842
$ws.Range("T9").Value = 0.002533870979031098
$ws.Range("U9").Value = 0.002470804824265213
$ws.Range("V9").Value = 0.002470804824265213
$ws.Range("W9").Value = 0.002456415562899114
$ws.Range("X9").Value = 0.002402388886296047
$ws.Range("Y9").Value = 0.002400768087259848
$ws.Range("C10").Value = 0.222224235534668
$ws.Range("E10").Value = 131.3187775684819
$ws.Range("F10").Value = 0.005332121010410525
$ws.Range("G10").Value = 0.004086079520537409
$ws.Range("H10").Value = 0.003705523309299093
$ws.Range("I10").Value = 0.003461786051366342
$ws.Range("J10").Value = 0.003419299326158093
$ws.Range("K10").Value = 0.003385313938731848
$ws.Range("L10").Value = 0.003057548374169295
$ws.Range("M10").Value = 0.003057548374169295
$ws.Range("N10").Value = 0.002881505085931732
$ws.Range("O10").Value = 0.002881505085931732
$ws.Range("P10").Value = 0.002809257586947448
$ws.Range("Q10").Value = 0.002729107386460642
$ws.Range("R10").Value = 0.002726573061178912
$ws.Range("S10").Value = 0.002726573061178912
$ws.Range("T10").Value = 0.002662052494544401
$ws.Range("U10").Value = 0.002652874126678084
$ws.Range("V10").Value = 0.002641424918753105
$ws.Range("W10").Value = 0.002586359042402178
$ws.Range("X10").Value = 0.002564638119029839
$ws.Range("Y10").Value = 0.002559820225506469
$ws.Range("C11").Value = 0.2078926563262939
$ws.Range("E11").Value = 130.6791473959038
$ws.Range("F11").Value = 0.004959683832992694
$ws.Range("G11").Value = 0.003888169743577927
$ws.Range("H11").Value = 0.003477358963056658
$ws.Range("I11").Value = 0.003477358963056658
$ws.Range("J11").Value = 0.003288669707095468
$ws.Range("K11").Value = 0.003014984416656826
$ws.Range("L11").Value = 0.002969485827017234
$ws.Range("M11").Value = 0.002969485827017234
$ws.Range("N11").Value = 0.002896639929711429
$ws.Range("O11").Value = 0.002804100637848087
$ws.Range("P11").Value = 0.002792435969355036
$ws.Range("Q11").Value = 0.002757108026358998
$ws.Range("R11").Value = 0.002741502051850018
$ws.Range("S11").Value = 0.002713232308939062
$ws.Range("T11").Value = 0.002710031841885551
$ws.Range("U11").Value = 0.002686023812329275
$ws.Range("V11").Value = 0.002640882957543881
$ws.Range("W11").Value = 0.002602610744273506
$ws.Range("X11").Value = 0.002580640908618788
$ws.Range("Y11").Value = 0.002547351801089742
